# Updates cryptos list data (prices + volume %) for rows 2-51,
# matching the GitHub Actions scheduled refresh commit.
# Numeric-looking price strings are force-written as text (NumberFormat "@")
# so trailing zeros / formatting are preserved exactly, then the cell style
# is reset to "Normal" so no stray explicit style index is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '97.052.12'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '3.708.36'
$ws.Range("E3").Value = '  +3.98%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.65%  '
$ws.Range("E6").Value = '  +18.76%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '660.75'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.429'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.09'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.65%  '
$ws.Range("E10").Value = '  -0.06%  '
$ws.Range("D11").Value = '3.712.16'
$ws.Range("E11").Value = '  +4.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.70'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.59%  '
$ws.Range("E13").Value = '  +1.54%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.71%  '
$ws.Range("D15").Value = '4.394.61'
$ws.Range("E15").Value = '  +3.88%  '
$ws.Range("D16").Value = '97.135.06'
$ws.Range("E16").Value = '  +1.49%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000262'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.10%  '
$ws.Range("D18").Value = '3.673.01'
$ws.Range("E18").Value = '  +3.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '18.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.02%  '
$ws.Range("E22").Value = '  +4.87%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.47'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.34%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '515.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000209'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.35%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '101.46'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +5.61%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '13.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.42%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.169'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.23%  '
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.77%  '
$ws.Range("E32").Value = '  +0.16%  '
$ws.Range("E33").Value = '  +2.19%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '33.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.74'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.78%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.591'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.72%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '619.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '43.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +27.67%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.159'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.975'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.95%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.21%  '
$ws.Range("E45").Value = '  +8.27%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0442'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.68%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.25%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.418'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +24.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.63'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.60'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.39%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.68'
$ws.Range("D51").Style = "Normal"
